# Generate Report for Handoff
# Replace the old handoff file identifier (6841b528-09e2-47df-86e5-8d173574fc02)
# with the new one (431d1843-bd14-45d9-9387-4015ca2b9a76) and its new xliff
# hashes/timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "6841b528-09e2-47df-86e5-8d173574fc02"
$newGuid = "431d1843-bd14-45d9-9387-4015ca2b9a76"

$oldHash = "f1926476c8bcdf58241fb074b42960d5adf86260"
$newHash = "400dcc228595f326ad3b27ed963e322bd1ab34a4"

$newFileName = $newGuid + ".md"
$newPathName = "e2e\" + $newGuid + ".md"

$newHoDate        = "2016-08-27 02:56:08"
$newZhCnXlf       = $newGuid + "." + $newHash + ".zh-cn.xlf"
$newZhCnHandoffDt = "2016-08-27 02:56:00"
$newDeDeXlf       = $newGuid + "." + $newHash + ".de-de.xlf"

# The hyperlink target URL itself is unchanged by this edit (it still points at
# the original commit path) - only the displayed text changes.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7efb8aa12b89ab668f3967caa320ea191f7b0183/e2e/" + $oldGuid + ".md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathName
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", $newPathName)

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newFileName
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnHandoffDt

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", $newFileName)

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newFileName
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newHoDate

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", $newFileName)
